$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5841821432113647
$ws.Range("B1").Value = 1.456153512001038
$ws.Range("C1").Value = 6.055192470550537
$ws.Range("D1").Value = 1.844111680984497
$ws.Range("E1").Value = 1.587704300880432
